# 9.5.2 indicator sheet update:
#  - replace the Russian header text in B1 with the new wording
#    (this also retires the old shared string and appends the new one,
#     matching the reference diff's shared-strings reordering)
#  - add a new "2023" data column (Q) with header style copied from P4
#    and data style copied from P5, then set the new values
#  - move the active selection back to A1 (closest available approximation
#    to the reference removing the saved selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Russian column header (B1) to the new wording.
$ws.Range("B1").Value = "9.5.2 Количество исследователей (в эквиваленте полной занятости) на миллион жителей"

# Add new column Q for year 2023, copying formatting from column P.
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2023

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 631

$excel.CutCopyMode = 0

[void]$ws.Range("A1").Select()
